$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "GuestList"

# Update selection
$ws.Range("G15").Select()

# Add new row of data
$ws.Range("A7").Value = "nisayon ert"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "bride"
